# Apply updated dSF (column F) values to reflect repulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -3
    4  = -4
    5  = -6
    7  = 3
    8  = -2
    9  = -8
    14 = 9
    15 = -3
    17 = 2
    18 = 1
    19 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
